# VfM class refactored to latest data structure
# The "A11 / HSMRPG" row (row 5) on the Q4_19_20 sheet is obsolete and is
# removed, which shifts the following rows (A13, F9, Columbia) up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Q4_19_20")

# Delete the entire row 5 (B5:C5 = "A11"/"HSMRPG"); this shifts rows 6-8
# upward to become rows 5-7 and updates the sheet's used range/dimension
# from B2:L8 to B2:L7 automatically.
$ws.Rows.Item(5).Delete()
